$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '61.572.69'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -3.66%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.001.98'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -2.55%  '

$ws.Range("E4").Value = '  +0.03%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '538.82'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.33%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '132.40'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -2.60%  '

$ws.Range("E7").Value = '  +0.10%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '2.997.68'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -2.56%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.494'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.12%  '

$ws.Range("E10").Value = '  -5.37%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '6.08'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -5.40%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.448'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -1.50%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.0000222'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -2.27%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '33.74'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -1.79%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '3.486.04'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -2.45%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '61.629.82'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -3.56%  '

$ws.Range("E17").Value = '  -2.20%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '2.999.86'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -2.61%  '

$ws.Range("E19").Value = '  -0.58%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '469.66'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -2.79%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '13.16'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -1.35%  '

$ws.Range("E22").Value = '  -4.10%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '6.95'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -1.85%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '80.25'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.66%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '11.96'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -1.32%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.05%  '

$ws.Range("E27").Value = '  -0.07%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '7.68'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -5.00%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.02%  '

$ws.Range("E30").Value = '  +0.54%  '

$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '25.54'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -2.55%  '

$ws.Range("B32").Value = 'Mantle'
$ws.Range("C32").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '1.15'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -1.02%  '

$ws.Range("B33").Value = 'Stacks'
$ws.Range("C33").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '2.29'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -2.20%  '

$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '55.42'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -2.90%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '5.35'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.18%  '

$ws.Range("E36").Value = '  -1.47%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '454.73'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -8.92%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '3.180.91'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -1.69%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.0791'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +0.15%  '

$ws.Range("E40").Value = '  +0.18%  '

$ws.Range("E41").Value = '  -4.00%  '

$ws.Range("E42").Value = '  -0.10%  '

$ws.Range("E43").Value = '  -9.47%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '25.47'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +4.76%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.241'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -4.90%  '

$ws.Range("B47").Value = 'Fetch.AI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '1.97'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -3.13%  '

$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '117.87'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -3.13%  '

$ws.Range("E49").Value = '  -1.03%  '

$ws.Range("E50").Value = '  -7.73%  '

$ws.Range("E51").Value = '  +5.69%  '
